$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2514
$ws.Range("J98").Value = 3496.5
$ws.Range("L98").Value = 3496.5
$ws.Range("N98").Value = -6492.5
$ws.Range("H122").Value = 2514
$ws.Range("J122").Value = 3496.5
$ws.Range("L122").Value = 10489.5
$ws.Range("N122").Value = -15389.5
$ws.Range("H127").Value = 528.5
$ws.Range("I127").Value = 528.5
$ws.Range("K127").Value = 1585.5
$ws.Range("M127").Value = 3374.5
$ws.Range("H138").Value = 12667.643
$ws.Range("J138").Value = 13174
$ws.Range("L138").Value = 39522
$ws.Range("N138").Value = -49802
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14925.154
$ws.Range("I32").Value = 6766.788
$ws.Range("J32").Value = 23338.469
$ws.Range("K32").Value = 6766.788
$ws.Range("L32").Value = 23338.469
$ws.Range("M32").Value = -6479.788
$ws.Range("N32").Value = -23912.469
$ws.Range("H122").Value = 387515.94
$ws.Range("I122").Value = 590272.9399999999
$ws.Range("K122").Value = 1770818.82
$ws.Range("M122").Value = -1768368.82
$ws.Range("H132").Value = 1827.2941
$ws.Range("I132").Value = 1783.84
$ws.Range("K132").Value = 5351.52
$ws.Range("M132").Value = -2821.52
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1125.5
$ws.Range("I86").Value = 1142.8462
$ws.Range("J86").Value = 900
$ws.Range("K86").Value = 1142.8462
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = -19.84619999999995
$ws.Range("N86").Value = -3146
$ws.Range("H89").Value = 1125.5
$ws.Range("I89").Value = 1142.8462
$ws.Range("J89").Value = 900
$ws.Range("K89").Value = 5714.231
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = -98.23099999999977
$ws.Range("N89").Value = -15732
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 18481.666
$ws.Range("I69").Value = 18481.666
$ws.Range("K69").Value = 18481.666
$ws.Range("M69").Value = -17732.666
$ws.Range("H72").Value = 18481.666
$ws.Range("I72").Value = 18481.666
$ws.Range("K72").Value = 55444.99800000001
$ws.Range("M72").Value = -51700.99800000001
$ws.Range("H99").Value = 15656.353
$ws.Range("I99").Value = 15096
$ws.Range("J99").Value = 16048.6
$ws.Range("K99").Value = 15096
$ws.Range("L99").Value = 16048.6
$ws.Range("M99").Value = -13598
$ws.Range("N99").Value = -19044.6
$ws.Range("H126").Value = 15656.353
$ws.Range("I126").Value = 15096
$ws.Range("J126").Value = 16048.6
$ws.Range("K126").Value = 45288
$ws.Range("L126").Value = 48145.8
$ws.Range("M126").Value = -42818
$ws.Range("N126").Value = -53085.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 399.75
$ws.Range("I5").Value = 399.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1199.25
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1087.25
$ws.Range("N5").ClearContents()
$ws.Range("H75").Value = 1024.4
$ws.Range("I75").Value = 649
$ws.Range("J75").Value = 1118.25
$ws.Range("K75").Value = 1947
$ws.Range("L75").Value = 3354.75
$ws.Range("M75").Value = -949
$ws.Range("N75").Value = -5350.75
$ws.Range("H78").Value = 1024.4
$ws.Range("I78").Value = 649
$ws.Range("J78").Value = 1118.25
$ws.Range("K78").Value = 5841
$ws.Range("L78").Value = 10064.25
$ws.Range("M78").Value = -849
$ws.Range("N78").Value = -20048.25
$ws.Range("H113").Value = 2014.1666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2014.1666
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6042.4998
$ws.Range("N113").Value = -10382.4998
$ws.Range("M113").ClearContents()
$ws.Range("H131").Value = 1559.8
$ws.Range("I131").Value = 622.8
$ws.Range("J131").Value = 2028.3
$ws.Range("K131").Value = 1868.4
$ws.Range("L131").Value = 6084.9
$ws.Range("M131").Value = 3171.6
$ws.Range("N131").Value = -16164.9
$ws.Range("H132").Value = 2270.8333
$ws.Range("J132").Value = 8333.333000000001
$ws.Range("L132").Value = 74999.997
$ws.Range("N132").Value = -80059.997
$ws.Range("H133").Value = 18343.334
$ws.Range("I133").Value = 15030
$ws.Range("J133").Value = 20000
$ws.Range("K133").Value = 45090
$ws.Range("L133").Value = 60000
$ws.Range("M133").Value = -40030
$ws.Range("N133").Value = -70120
$ws.Range("H135").Value = 399.75
$ws.Range("I135").Value = 399.75
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3597.75
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1062.75
$ws.Range("N135").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 552474.4
$ws.Range("I122").Value = 78869.16
$ws.Range("K122").Value = 236607.48
$ws.Range("M122").Value = -234157.48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4231.4443
$ws.Range("J7").Value = 5496.6665
$ws.Range("L7").Value = 5496.6665
$ws.Range("N7").Value = -5720.6665
$ws.Range("H16").Value = 6141.533
$ws.Range("I16").Value = 6438.0713
$ws.Range("J16").Value = 1990
$ws.Range("K16").Value = 6438.0713
$ws.Range("L16").Value = 1990
$ws.Range("M16").Value = -6268.0713
$ws.Range("N16").Value = -2330
$ws.Range("H22").Value = 917.2857
$ws.Range("I22").Value = 972
$ws.Range("K22").Value = 972
$ws.Range("M22").Value = -677
$ws.Range("H27").Value = 917.2857
$ws.Range("I27").Value = 972
$ws.Range("K27").Value = 972
$ws.Range("M27").Value = -865
$ws.Range("H68").Value = 3699.8
$ws.Range("I68").Value = 3499
$ws.Range("J68").Value = 3750
$ws.Range("K68").Value = 3499
$ws.Range("L68").Value = 3750
$ws.Range("M68").Value = -2750
$ws.Range("N68").Value = -5248
$ws.Range("H71").Value = 3699.8
$ws.Range("I71").Value = 3499
$ws.Range("J71").Value = 3750
$ws.Range("K71").Value = 17495
$ws.Range("L71").Value = 18750
$ws.Range("M71").Value = -13751
$ws.Range("N71").Value = -26238
$ws.Range("H126").Value = 4231.4443
$ws.Range("J126").Value = 5496.6665
$ws.Range("L126").Value = 16489.9995
$ws.Range("N126").Value = -21429.9995
$ws.Range("H136").Value = 3737.9
$ws.Range("I136").Value = 3597.6667
$ws.Range("K136").Value = 10793.0001
$ws.Range("M136").Value = -8243.000100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H49").Value = 28000
$ws.Range("J49").Value = 28000
$ws.Range("L49").Value = 28000
$ws.Range("N49").Value = -28460
$ws.Range("H119").Value = 40749.75
$ws.Range("J119").Value = 40749.75
$ws.Range("L119").Value = 40749.75
$ws.Range("N119").Value = -50425.75
